# Update the recomputed power-flow results (pl_mw) for the 380 kV case.
# Only cells B, C, D, E, F, H, L, N change for rows 2-25 (A = 0..23),
# matching the "case with 380 kV done" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 0)
$ws.Range("B2").Value = 2.841834217422502
$ws.Range("C2").Value = 0.2377829919085457
$ws.Range("D2").Value = 0.1132707645347892
$ws.Range("E2").Value = 0.049046120412779
$ws.Range("F2").Value = 2.421976638479862
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("L2").Value = 0.2716327507868641
$ws.Range("N2").Value = 1.946590460329162

# Row 3 (A3 = 1)
$ws.Range("B3").Value = 2.683083395745541
$ws.Range("C3").Value = 0.2069324353714705
$ws.Range("D3").Value = 0.1136471591582904
$ws.Range("E3").Value = 0.04924918981984883
$ws.Range("F3").Value = 2.364479096077304
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("L3").Value = 0.2616684373783045
$ws.Range("N3").Value = 1.958222973137481

# Row 4 (A4 = 2)
$ws.Range("B4").Value = 2.587462152768467
$ws.Range("C4").Value = 0.1880198126438586
$ws.Range("D4").Value = 0.1139298444661492
$ws.Range("E4").Value = 0.04938467912778588
$ws.Range("F4").Value = 2.330742348702188
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("L4").Value = 0.255727701541403
$ws.Range("N4").Value = 1.966024942044065

# Row 5 (A5 = 3)
$ws.Range("B5").Value = 2.548957376177782
$ws.Range("C5").Value = 0.1803193773987459
$ws.Range("D5").Value = 0.1140579358836504
$ws.Range("E5").Value = 0.04944261196032196
$ws.Range("F5").Value = 2.317385159484374
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("L5").Value = 0.253351109369504
$ws.Range("N5").Value = 1.969369320993898

# Row 6 (A6 = 4)
$ws.Range("B6").Value = 2.542591450739678
$ws.Range("C6").Value = 0.1790411002427561
$ws.Range("D6").Value = 0.114079982016861
$ws.Range("E6").Value = 0.04945239602782303
$ws.Range("F6").Value = 2.315190718983146
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("L6").Value = 0.2529591453428708
$ws.Range("N6").Value = 1.969934594587173

# Row 7 (A7 = 5)
$ws.Range("B7").Value = 2.586941000068123
$ws.Range("C7").Value = 0.187915936247947
$ws.Range("D7").Value = 0.1139315198387649
$ws.Range("E7").Value = 0.04938544941335987
$ws.Range("F7").Value = 2.330560630820443
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("L7").Value = 0.2556954709875185
$ws.Range("N7").Value = 1.966069378471538

# Row 8 (A8 = 6)
$ws.Range("B8").Value = 2.786709895397905
$ws.Range("C8").Value = 0.227139012690742
$ws.Range("D8").Value = 0.1133897951206606
$ws.Range("E8").Value = 0.04911389825120471
$ws.Range("F8").Value = 2.401824287529735
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("L8").Value = 0.2681600511346431
$ws.Range("N8").Value = 1.950463952926299

# Row 9 (A9 = 7)
$ws.Range("B9").Value = 3.193367121953372
$ws.Range("C9").Value = 0.304334464293845
$ws.Range("D9").Value = 0.1127401793197365
$ws.Range("E9").Value = 0.04866697884750781
$ws.Range("F9").Value = 2.554161799494267
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("L9").Value = 0.2940259638491654
$ws.Range("N9").Value = 1.925133008421824

# Row 10 (A10 = 8)
$ws.Range("B10").Value = 3.501561083298441
$ws.Range("C10").Value = 0.3612856057851559
$ws.Range("D10").Value = 0.1125194272341474
$ws.Range("E10").Value = 0.04839063558229384
$ws.Range("F10").Value = 2.673992092429216
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("L10").Value = 0.3139209947141666
$ws.Range("N10").Value = 1.909788611674941

# Row 11 (A11 = 9)
$ws.Range("B11").Value = 3.643888120940574
$ws.Range("C11").Value = 0.3872600052347366
$ws.Range("D11").Value = 0.1124757954356213
$ws.Range("E11").Value = 0.04827617917381932
$ws.Range("F11").Value = 2.730275165508971
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("L11").Value = 0.3231707290510428
$ws.Range("N11").Value = 1.90352888131163

# Row 12 (A12 = 10)
$ws.Range("B12").Value = 3.698095256445981
$ws.Range("C12").Value = 0.397106575249154
$ws.Range("D12").Value = 0.1124675299819984
$ws.Range("E12").Value = 0.04823445324871578
$ws.Range("F12").Value = 2.751846892124007
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("L12").Value = 0.3267024377777261
$ws.Range("N12").Value = 1.901263075431615

# Row 13 (A13 = 11)
$ws.Range("B13").Value = 3.686406868920017
$ws.Range("C13").Value = 0.3949854515706193
$ws.Range("D13").Value = 0.1124689414019855
$ws.Range("E13").Value = 0.04824336782186922
$ws.Range("F13").Value = 2.747189474832993
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("L13").Value = 0.3259405238874393
$ws.Range("N13").Value = 1.901746388098829

# Row 14 (A14 = 12)
$ws.Range("B14").Value = 3.648341513346054
$ws.Range("C14").Value = 0.3880698697322487
$ws.Range("D14").Value = 0.1124749494244881
$ws.Range("E14").Value = 0.04827271397643607
$ws.Range("F14").Value = 2.732044680831706
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("L14").Value = 0.3234607008323707
$ws.Range("N14").Value = 1.903340369560397

# Row 15 (A15 = 13)
$ws.Range("B15").Value = 3.625066031417191
$ws.Range("C15").Value = 0.3838352899968527
$ws.Range("D15").Value = 0.1124797075340211
$ws.Range("E15").Value = 0.04829089976659251
$ws.Range("F15").Value = 2.722801841294512
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("L15").Value = 0.3219455297762863
$ws.Range("N15").Value = 1.90433038381741

# Row 16 (A16 = 14)
$ws.Range("B16").Value = 3.492302920214684
$ws.Range("C16").Value = 0.3595895507752971
$ws.Range("D16").Value = 0.1125234291385766
$ws.Range("E16").Value = 0.04839834185333114
$ws.Range("F16").Value = 2.670349833742421
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("L16").Value = 0.3133205458588435
$ws.Range("N16").Value = 1.910212284256033

# Row 17 (A17 = 15)
$ws.Range("B17").Value = 3.411405592728556
$ws.Range("C17").Value = 0.3447334429776561
$ws.Range("D17").Value = 0.1125648600629177
$ws.Range("E17").Value = 0.04846713492839871
$ws.Range("F17").Value = 2.638628664890916
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("L17").Value = 0.3080807210012324
$ws.Range("N17").Value = 1.914005893466936

# Row 18 (A18 = 16)
$ws.Range("B18").Value = 3.365075588672994
$ws.Range("C18").Value = 0.3361948525414959
$ws.Range("D18").Value = 0.1125940297986148
$ws.Range("E18").Value = 0.0485077622351282
$ws.Range("F18").Value = 2.620550015391245
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("L18").Value = 0.3050856571229161
$ws.Range("N18").Value = 1.916255633949632

# Row 19 (A19 = 17)
$ws.Range("B19").Value = 3.349423246077208
$ws.Range("C19").Value = 0.333304876140744
$ws.Range("D19").Value = 0.112604820543531
$ws.Range("E19").Value = 0.04852169995042832
$ws.Range("F19").Value = 2.614457382392175
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("L19").Value = 0.3040747884484318
$ws.Range("N19").Value = 1.917028964194643

# Row 20 (A20 = 18)
$ws.Range("B20").Value = 3.419996516625929
$ws.Range("C20").Value = 0.3463142470774301
$ws.Range("D20").Value = 0.1125598963739449
$ws.Range("E20").Value = 0.04845970216427986
$ws.Range("F20").Value = 2.641988175277419
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("L20").Value = 0.3086365666805051
$ws.Range("N20").Value = 1.913595037459942

# Row 21 (A21 = 19)
$ws.Range("B21").Value = 3.659513751074769
$ws.Range("C21").Value = 0.3901008473259822
$ws.Range("D21").Value = 0.1124729598992076
$ws.Range("E21").Value = 0.04826405045083071
$ws.Range("F21").Value = 2.736486026520964
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("L21").Value = 0.3241882940303498
$ws.Range("N21").Value = 1.90286933113623

# Row 22 (A22 = 20)
$ws.Range("B22").Value = 3.817867099816681
$ws.Range("C22").Value = 0.4187803459709585
$ws.Range("D22").Value = 0.1124643131141312
$ws.Range("E22").Value = 0.04814560052720918
$ws.Range("F22").Value = 2.799754657441071
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("L22").Value = 0.3345216605275851
$ws.Range("N22").Value = 1.896469736218307

# Row 23 (A23 = 21)
$ws.Range("B23").Value = 3.733183194872481
$ws.Range("C23").Value = 0.4034675144938547
$ws.Range("D23").Value = 0.1124644899724316
$ws.Range("E23").Value = 0.04820795820745172
$ws.Range("F23").Value = 2.765847672990844
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("L23").Value = 0.3289909245548301
$ws.Range("N23").Value = 1.899829142756928

# Row 24 (A24 = 22)
$ws.Range("B24").Value = 3.416112002925047
$ws.Range("C24").Value = 0.3455995581934417
$ws.Range("D24").Value = 0.1125621237983339
$ws.Range("E24").Value = 0.04846305916142857
$ws.Range("F24").Value = 2.6404688479318
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("L24").Value = 0.3083852147374841
$ws.Range("N24").Value = 1.913780571414875

# Row 25 (A25 = 23)
$ws.Range("B25").Value = 3.081725240779633
$ws.Range("C25").Value = 0.2834144603444315
$ws.Range("D25").Value = 0.1128712313768077
$ws.Range("E25").Value = 0.04877873720234671
$ws.Range("F25").Value = 2.511578800826101
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("L25").Value = 0.2868737749718377
$ws.Range("N25").Value = 1.931415975659021
